$wb = $excel.ActiveWorkbook

# --- Asset_Cal_Info: corrected FLORT calibration values ---
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_scattering_angle (row 7) corrected from 117 to 124
$assetCal.Range("F7").Value = 124

# CC_angular_resolution (row 9) corrected from 1.08 to 1.076
$assetCal.Range("F9").Value = 1.076

# --- Make Asset_Cal_Info the active/selected sheet & set its selection ---
$assetCal.Activate()
$assetCal.Range("E25").Select()
